$wb = $excel.ActiveWorkbook

# Update the "展览" sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1287
$ws1.Range("F3").Value = 1650
$ws1.Range("F5").Value = 6214
$ws1.Range("F6").Value = 43

# Update the "全部类型" sheet (mirrors the same data)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1287
$ws4.Range("F3").Value = 1650
$ws4.Range("F5").Value = 6214
$ws4.Range("F6").Value = 43
